$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '54.247.73'
$ws.Range('E2').Value = '  -2.66%  '
$ws.Range('D3').Value = '2.286.54'
$ws.Range('E3').Value = '  -2.18%  '
$ws.Range('D4').Value = "'0.999"
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').Value = "'495.23"
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -1.41%  '
$ws.Range('D6').Value = "'127.51"
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -0.76%  '
$ws.Range('D7').Value = "'0.998"
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('E8').Value = '  -1.22%  '
$ws.Range('D9').Value = '2.283.47'
$ws.Range('E9').Value = '  -2.48%  '
$ws.Range('E10').Value = '  -3.87%  '
$ws.Range('E11').Value = '  +0.21%  '
$ws.Range('E12').Value = '  +0.68%  '
$ws.Range('D13').Value = "'4.65"
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -3.43%  '
$ws.Range('D14').Value = '2.686.90'
$ws.Range('E14').Value = '  -2.38%  '
$ws.Range('D15').Value = "'21.63"
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +0.33%  '
$ws.Range('D16').Value = '54.142.96'
$ws.Range('E16').Value = '  -2.79%  '
$ws.Range('E17').Value = '  -2.10%  '
$ws.Range('D18').Value = '2.283.66'
$ws.Range('E18').Value = '  +0.09%  '
$ws.Range('D19').Value = "'9.92"
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -0.13%  '
$ws.Range('D20').Value = "'4.03"
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +0.76%  '
$ws.Range('D21').Value = "'298.84"
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -3.37%  '
$ws.Range('D22').Value = "'6.28"
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +1.37%  '
$ws.Range('E23').Value = '  -0.06%  '
$ws.Range('D24').Value = "'63.81"
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -2.09%  '
$ws.Range('E25').Value = '  +0.66%  '
$ws.Range('E26').Value = '  +1.40%  '
$ws.Range('D27').Value = '2.390.18'
$ws.Range('E27').Value = '  -2.30%  '
$ws.Range('E28').Value = '  +1.46%  '
$ws.Range('D29').Value = "'7.12"
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +0.12%  '
$ws.Range('D30').Value = "'163.48"
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -5.12%  '
$ws.Range('D31').Value = "'1.61"
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -1.78%  '
$ws.Range('E32').Value = '  -2.66%  '
$ws.Range('D33').Value = "'5.85"
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +1.28%  '
$ws.Range('D34').Value = "'0.999"
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +0.00%  '
$ws.Range('D35').Value = "'0.998"
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +0.15%  '
$ws.Range('D36').Value = "'1.06"
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +0.70%  '
$ws.Range('D37').Value = "'17.48"
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -0.60%  '
$ws.Range('E38').Value = '  +0.49%  '
$ws.Range('D39').Value = "'0.869"
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +5.22%  '
$ws.Range('E40').Value = '  -0.35%  '
$ws.Range('D41').Value = "'35.33"
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -2.12%  '
$ws.Range('E42').Value = '  +1.70%  '
$ws.Range('E43').Value = '  +1.28%  '
$ws.Range('E44').Value = '  -0.36%  '
$ws.Range('B45').Value = 'Aave'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D45').Value = "'126.30"
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -0.60%  '
$ws.Range('B46').Value = 'RenderToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D46').Value = "'4.85"
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +3.39%  '
$ws.Range('D47').Value = "'0.0889"
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -0.39%  '
$ws.Range('D48').Value = "'0.547"
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -1.45%  '
$ws.Range('D49').Value = "'238.56"
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +0.86%  '
$ws.Range('E50').Value = '  +0.60%  '
$ws.Range('E51').Value = '  -0.88%  '
